# Insert a new row of weekly price data at row 344, pushing the existing
# rows 344-354 down to 345-355.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 344 (existing rows shift down).
$ws.Rows(344).Insert()

# Copy the date-column number format/style from the row below (old row 344,
# now row 345) so the new date cell renders the same way.
$ws.Cells.Item(345, 4).Copy()
$ws.Cells.Item(344, 4).PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row with the new weekly record.
$ws.Cells.Item(344, 1).Value = 10
$ws.Cells.Item(344, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(344, 3).Value = "La Araucanía"
$ws.Cells.Item(344, 4).Value = 45239
$ws.Cells.Item(344, 5).Value = 9
$ws.Cells.Item(344, 6).Value = 100112005
$ws.Cells.Item(344, 7).Value = "Puerro"
$ws.Cells.Item(344, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(344, 9).Value = "Primera"
$ws.Cells.Item(344, 10).Value = 90
$ws.Cells.Item(344, 11).Value = 10000
$ws.Cells.Item(344, 12).Value = 10000
$ws.Cells.Item(344, 13).Value = 10000
$ws.Cells.Item(344, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(344, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(344, 16).Value = 833
$ws.Cells.Item(344, 17).Value = 12
$ws.Cells.Item(344, 18).Value = "Hortaliza"
